$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 17:38"

# Full refreshed country table (rank-sorted by "Casos totales" desc), rows 4.. onward
$countries = @(
  @("Estados Unidos", 2904003, 13415, 1236019, 1535837, 0, 46, 132147),
  @("Brasil", 1545458, 2117, 978615, 503548, 0, 41, 63295),
  @("Rusia", 674515, 6632, 446879, 217609, 0, 168, 10027),
  @("India", 661233, 11344, 401847, 240474, 0, 243, 18912),
  @("España", 297625, 0, 0, 0, 0, 0, 28385),
  @("Peru", 295599, 0, 185852, 99521, 0, 0, 10226),
  @("Chile", 288089, 0, 253343, 28695, 0, 0, 6051),
  @("Reino Unido", 284276, 0, 0, 0, 0, 0, 44131),
  @("Mexico", 245251, 6740, 147205, 68203, 0, 654, 29843),
  @("Italia", 241184, 0, 191467, 14884, 0, 0, 34833),
  @("Iran", 237878, 2449, 198949, 27521, 0, 148, 11408),
  @("Pakistan", 225283, 3387, 125094, 95570, 0, 68, 4619),
  @("Arabia Saudita", 205929, 4128, 143256, 60815, 0, 56, 1858),
  @("Turquia", 203456, 0, 178278, 19992, 0, 0, 5186),
  @("Alemania", 197210, 210, 181300, 6836, 0, 1, 9074),
  @("Sudafrica", 177124, 0, 86298, 87874, 0, 0, 2952),
  @("Francia", 166960, 0, 77060, 60007, 0, 0, 29893),
  @("Banglades", 159679, 3288, 70721, 86961, 0, 29, 1997),
  @("Colombia", 109505, 0, 45334, 60394, 0, 0, 3777),
  @("Canada", 105091, 0, 68693, 27735, 0, 0, 8663),
  @("Catar", 99183, 530, 90387, 8673, 0, 2, 123),
  @("China", 83545, 3, 78509, 402, 0, 0, 4634),
  @("Argentina", 72786, 0, 25930, 45403, 0, 16, 1453),
  @("Egipto", 72711, 0, 19690, 49820, 0, 0, 3201),
  @("Suecia", 71419, 0, 0, 0, 0, 0, 5420),
  @("Bielorrusia", 63270, 273, 50669, 12183, 0, 6, 418),
  @("Indonesia", 62142, 1447, 28219, 30834, 0, 53, 3089),
  @("Belgica", 61838, 111, 17091, 34976, 0, 6, 9771),
  @("Ecuador", 60657, 0, 28391, 27566, 0, 0, 4700),
  @("Irak", 58354, 2334, 31077, 24909, 0, 106, 2368),
  @("Emiratos Arabes Unidos", 50857, 716, 39857, 10679, 0, 3, 321),
  @("Paises Bajos", 50335, 0, 0, 0, 0, 0, 6113),
  @("Kuwait", 49303, 631, 39943, 8995, 0, 5, 365),
  @("Ucrania", 47677, 914, 21115, 25335, 0, 15, 1227),
  @("Kazajistan", 45719, 1644, 26750, 18781, 0, 0, 188),
  @("Oman", 45106, 1177, 26968, 17935, 0, 10, 203),
  @("Singapur", 44664, 185, 39769, 4869, 0, 0, 26),
  @("Portugal", 43569, 413, 28772, 13192, 0, 7, 1605),
  @("Filipinas", 41830, 1494, 11453, 29087, 0, 10, 1290),
  @("Bolivia", 36818, 1290, 10766, 24732, 0, 49, 1320),
  @("Republica Dominicana", 36184, 1036, 18602, 16796, 0, 11, 786),
  @("Panama", 35995, 0, 16945, 18352, 0, 0, 698),
  @("Polonia", 35719, 314, 23127, 11080, 0, 5, 1512),
  @("Afganistan", 32672, 348, 19164, 12682, 0, 7, 826),
  @("Suiza", 32198, 97, 29200, 1033, 0, 0, 1965),
  @("Rumania", 28582, 416, 19854, 6997, 0, 23, 1731),
  @("Barein", 28410, 0, 23318, 4997, 0, 0, 95),
  @("Israel", 28055, 0, 17669, 10060, 0, 0, 326),
  @("Armenia", 27900, 580, 15935, 11488, 0, 8, 477),
  @("Nigeria", 27564, 0, 11069, 15867, 0, 0, 628),
  @("Irlanda", 25498, 0, 23364, 394, 0, 0, 1740),
  @("Honduras", 22116, 996, 2250, 19261, 0, 14, 605),
  @("Guatemala", 21293, 0, 3315, 17098, 0, 0, 880),
  @("Ghana", 19388, 0, 14330, 4941, 0, 0, 117),
  @("Azerbaiyan", 19267, 0, 10820, 8212, 0, 0, 235),
  @("Japon", 19068, 0, 16866, 1226, 0, 0, 976),
  @("Austria", 18165, 115, 16607, 853, 0, 0, 705),
  @("Moldavia", 17672, 227, 10396, 6696, 0, 8, 580),
  @("Serbia", 15829, 325, 13064, 2459, 0, 8, 306),
  @("Nepal", 15491, 232, 6415, 9042, 0, 2, 34),
  @("Argelia", 15070, 0, 10832, 3301, 0, 0, 937),
  @("Marruecos", 13434, 146, 9214, 3988, 0, 2, 232),
  @("Corea del Sur", 13030, 63, 11811, 936, 0, 1, 283),
  @("Dinamarca", 12832, 0, 11817, 409, 0, 0, 606),
  @("Camerun", 12592, 0, 10100, 2179, 0, 0, 313),
  @("Chequia", 12372, 53, 7851, 4168, 0, 0, 353),
  @("Costa de Marfil", 10244, 0, 4726, 5448, 0, 0, 70),
  @("Sudan", 9663, 0, 4624, 4435, 0, 0, 604),
  @("Uzbekistan", 9500, 104, 6251, 3220, 0, 0, 29),
  @("Noruega", 8921, 0, 8138, 532, 0, 0, 251),
  @("Malasia", 8658, 10, 8461, 76, 0, 0, 121),
  @("Australia", 8362, 107, 7355, 903, 0, 0, 104),
  @("Kenia", 7577, 389, 2236, 5182, 0, 5, 159),
  @("El Salvador", 7507, 240, 4434, 2863, 0, 8, 210),
  @("Consejo Danes para los Refugiados", 7379, 68, 2961, 4236, 0, 3, 182),
  @("Senegal", 7272, 108, 4713, 2430, 0, 4, 129),
  @("Finlandia", 7248, 6, 6700, 219, 0, 0, 329),
  @("Kirguistan", 7094, 327, 2714, 4302, 0, 2, 78),
  @("Republica de Macedonia", 6787, 0, 2876, 3583, 0, 0, 328),
  @("Venezuela", 6537, 0, 2100, 4378, 0, 0, 59),
  @("Haiti", 6230, 129, 1286, 4834, 0, 0, 110),
  @("Tayikistan", 6058, 0, 4690, 1316, 0, 0, 52),
  @("Etiopia", 5846, 0, 2430, 3313, 0, 0, 103),
  @("Gabon", 5620, 0, 2555, 3021, 0, 0, 44),
  @("Guinea", 5521, 0, 4446, 1042, 0, 0, 33),
  @("Bulgaria", 5497, 0, 2892, 2366, 0, 0, 239),
  @("Bosnia y Herzegovina", 4962, 0, 2550, 2221, 0, 0, 191),
  @("Republica de Yibuti", 4736, 0, 4580, 101, 0, 0, 55),
  @("Mauritania", 4705, 0, 1765, 2811, 0, 0, 129),
  @("Guayana Francesa", 4558, 0, 1777, 2765, 0, 0, 16),
  @("Luxemburgo", 4476, 29, 4016, 350, 0, 0, 110),
  @("Costa Rica", 4311, 0, 1657, 2636, 0, 0, 18),
  @("Hungria", 4174, 2, 2784, 801, 0, 1, 589),
  @("Republica de Africa Central", 3918, 0, 859, 3011, 0, 0, 48),
  @("Estado de Palestina", 3599, 265, 463, 3124, 0, 1, 12),
  @("Grecia", 3511, 25, 1374, 1945, 0, 0, 192),
  @("Tailandia", 3185, 5, 3066, 61, 0, 0, 58),
  @("Croacia", 3094, 86, 2183, 798, 0, 1, 113),
  @("Guinea Ecuatorial", 3071, 0, 842, 2178, 0, 0, 51),
  @("Somalia", 2944, 0, 951, 1903, 0, 0, 90),
  @("Albania", 2819, 67, 1637, 1108, 0, 2, 74),
  @("Madagascar", 2728, 216, 1078, 1621, 0, 3, 29),
  @("Mayotte", 2661, 0, 2375, 251, 0, 0, 35),
  @("Nicaragua", 2519, 0, 1238, 1198, 0, 0, 83),
  @("Maldivas", 2410, 0, 1976, 424, 0, 0, 10),
  @("Cuba", 2369, 8, 2227, 56, 0, 0, 86),
  @("Paraguay", 2349, 0, 1113, 1217, 0, 0, 19),
  @("Mali", 2285, 0, 1507, 661, 0, 0, 117),
  @("Sri Lanka", 2074, 5, 1885, 178, 0, 0, 11),
  @("Sudan del Sur", 2021, 0, 333, 1650, 0, 0, 38),
  @("Estonia", 1993, 2, 1870, 54, 0, 0, 69),
  @("Libano", 1855, 25, 1304, 516, 0, 0, 35),
  @("Islandia", 1855, 0, 1832, 13, 0, 0, 10),
  @("Lituania", 1831, 3, 1545, 207, 0, 0, 79),
  @("Guinea-Bisau", 1765, 0, 676, 1064, 0, 0, 25),
  @("Eslovaquia", 1749, 29, 1466, 255, 0, 0, 28),
  @("Eslovenia", 1679, 29, 1384, 184, 0, 0, 111),
  @("Zambia", 1632, 0, 1348, 254, 0, 0, 30),
  @("Congo", 1557, 0, 501, 1012, 0, 0, 44),
  @("Nueva Zelanda", 1530, 0, 1490, 18, 0, 0, 22),
  @("Sierra Leona", 1524, 0, 1042, 420, 0, 0, 62),
  @("Malaui", 1498, 0, 317, 1165, 0, 0, 16),
  @("Cabo Verde", 1382, 0, 629, 738, 0, 0, 15),
  @("Hong Kong", 1259, 11, 1145, 107, 0, 0, 7),
  @("Yemen", 1240, 0, 536, 369, 0, 0, 335),
  @("Benin", 1199, 0, 333, 845, 0, 0, 21),
  @("Tunez", 1181, 0, 1045, 86, 0, 0, 50),
  @("Jordania", 1147, 0, 897, 240, 0, 0, 10),
  @("Letonia", 1123, 1, 1000, 93, 0, 0, 30),
  @("Niger", 1082, 0, 963, 51, 0, 0, 68),
  @("Ruanda", 1081, 0, 512, 566, 0, 0, 3),
  @("Republica de Chipre", 999, 0, 839, 141, 0, 0, 19),
  @("Burkina Faso", 980, 0, 852, 75, 0, 0, 53),
  @("Uruguay", 952, 0, 837, 87, 0, 0, 28),
  @("Georgia", 948, 5, 825, 108, 0, 0, 15),
  @("Mozambique", 939, 0, 249, 684, 0, 0, 6),
  @("Uganda", 927, 16, 868, 59, 0, 0, 0),
  @("Libia", 918, 0, 230, 661, 0, 0, 27),
  @("Suazilandia", 909, 0, 515, 381, 0, 0, 13),
  @("Republica del Chad", 871, 0, 786, 11, 0, 0, 74),
  @("Liberia", 869, 36, 369, 463, 0, 0, 37),
  @("Principado de Andorra", 855, 0, 800, 3, 0, 0, 52),
  @("Jamaica", 721, 6, 565, 146, 0, 0, 10),
  @("Santo Tome y Principe", 719, 0, 267, 439, 0, 0, 13),
  @("Crucero", 712, 0, 651, 48, 0, 0, 13),
  @("San Marino", 698, 0, 656, 0, 0, 0, 42),
  @("Malta", 672, 0, 650, 13, 0, 0, 9),
  @("Togo", 671, 0, 424, 233, 0, 0, 14),
  @("Montenegro", 663, 0, 315, 335, 0, 0, 13),
  @("Zimbabue", 625, 0, 176, 442, 0, 0, 7),
  @("Surinam", 561, 0, 267, 281, 0, 0, 13),
  @("Reunion", 533, 0, 472, 59, 0, 0, 2),
  @("Tanzania", 509, 0, 183, 305, 0, 0, 21),
  @("Taiwan", 449, 0, 438, 4, 0, 0, 7),
  @("Namibia", 375, 25, 25, 350, 0, 0, 0),
  @("Vietnam", 355, 0, 340, 15, 0, 0, 0),
  @("Mauricio", 341, 0, 330, 1, 0, 0, 10),
  @("Siria", 338, 10, 123, 205, 0, 0, 10),
  @("Isla de Man", 336, 0, 312, 0, 0, 0, 24),
  @("Angola", 328, 0, 107, 203, 0, 0, 18),
  @("Birmania", 313, 7, 240, 67, 0, 0, 6),
  @("Comoras", 309, 0, 241, 61, 0, 0, 7),
  @("Botsuana", 277, 0, 29, 247, 0, 0, 1),
  @("Guyana", 256, 0, 117, 125, 0, 0, 14),
  @("Martinica", 249, 0, 98, 137, 0, 0, 14),
  @("Mongolia", 220, 0, 183, 37, 0, 0, 0),
  @("Eritrea", 215, 0, 56, 159, 0, 0, 0),
  @("Islas Caimanes", 201, 0, 194, 6, 0, 0, 1),
  @("Burundi", 191, 0, 118, 72, 0, 0, 1),
  @("Islas Feroe", 187, 0, 187, 0, 0, 0, 0),
  @("Guadalupe", 184, 0, 157, 13, 0, 0, 14),
  @("Gibraltar", 178, 0, 176, 2, 0, 0, 0),
  @("Bermudas", 146, 0, 137, 0, 0, 0, 9),
  @("Camboya", 141, 0, 131, 10, 0, 0, 0),
  @("Brunei", 141, 0, 138, 0, 0, 0, 3),
  @("Trinidad yTobago", 130, 0, 115, 7, 0, 0, 8),
  @("Monaco", 106, 0, 95, 7, 0, 0, 4),
  @("Bahamas", 104, 0, 89, 4, 0, 0, 11),
  @("Aruba", 104, 0, 98, 3, 0, 0, 3),
  @("Barbados", 97, 0, 90, 0, 0, 0, 7),
  @("Liechtenstein", 83, 0, 81, 1, 0, 0, 1),
  @("Seychelles", 81, 0, 11, 70, 0, 0, 0),
  @("Butan", 78, 1, 51, 27, 0, 0, 0),
  @("San Martin (Parte Holandesa)", 78, 1, 63, 0, 0, 0, 15),
  @("Antigua y Barbuda", 68, 0, 23, 42, 0, 0, 3),
  @("Polinesia Francesa", 62, 0, 60, 2, 0, 0, 0),
  @("Gambia", 55, 0, 27, 26, 0, 0, 2),
  @("Macao", 46, 0, 45, 1, 0, 0, 0),
  @("Islas Turcas y Caicos", 45, 1, 11, 32, 0, 0, 2),
  @("San Martin (Parte Francesa)", 43, 0, 37, 3, 0, 0, 3),
  @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
  @("Lesoto", 35, 0, 11, 24, 0, 0, 0),
  @("Guam", 32, 0, 0, 31, 0, 0, 1),
  @("Belice", 30, 2, 19, 9, 0, 0, 2),
  @("San Vicente y las Granadinas", 29, 0, 29, 0, 0, 0, 0),
  @("Timor Oriental", 24, 0, 24, 0, 0, 0, 0),
  @("Curazao", 23, 0, 19, 3, 0, 0, 1),
  @("Granada", 23, 0, 23, 0, 0, 0, 0),
  @("Santa Lucia", 22, 0, 19, 3, 0, 0, 0),
  @("Nueva Caledonia", 21, 0, 21, 0, 0, 0, 0),
  @("Laos", 19, 0, 19, 0, 0, 0, 0),
  @("Fiyi", 18, 0, 18, 0, 0, 0, 0),
  @("Dominica", 18, 0, 18, 0, 0, 0, 0),
  @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
  @("San Cristobal y Nieves", 15, 0, 15, 0, 0, 0, 0),
  @("Groenlandia", 13, 0, 13, 0, 0, 0, 0),
  @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0),
  @("Santa Sede", 12, 0, 12, 0, 0, 0, 0),
  @("Papua Nueva Guinea", 11, 0, 8, 3, 0, 0, 0),
  @("Montserrat", 11, 0, 10, 0, 0, 0, 1),
  @("Sahara Occidental", 10, 0, 8, 1, 0, 0, 1),
  @("Islas Virgenes Britanicas", 8, 0, 7, 0, 0, 0, 1),
  @("Bonaire, San Eustaquio y Saba", 7, 0, 7, 0, 0, 0, 0),
  @("San Bartolome", 6, 0, 6, 0, 0, 0, 0),
  @("Anguila", 3, 0, 3, 0, 0, 0, 0),
  @("San Pedro y Miquelon", 1, 0, 1, 0, 0, 0, 0)
)

$r = 4
foreach ($row in $countries) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
